$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44706
$ws.Range("L2").Value = 'Especial'
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 16000
$ws.Range("O2").Value = 16000
$ws.Range("P2").Value = 16000
$ws.Range("Q2").Value = '$/caja 18 kilos granel'
$ws.Range("R2").Value = 'Región de O''Higgins'
$ws.Range("S2").Value = 889
$ws.Range("T2").Value = 18

$ws.Range("D3").Value = 44706
$ws.Range("N3").Value = 12500
$ws.Range("O3").Value = 12500
$ws.Range("P3").Value = 12500
$ws.Range("Q3").Value = '$/caja 18 kilos granel'
$ws.Range("R3").Value = 'Región de O''Higgins'
$ws.Range("S3").Value = 694
$ws.Range("T3").Value = 18

$ws.Range("D4").Value = 44698
$ws.Range("L4").Value = 'Especial'
$ws.Range("N4").Value = 18000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 18000
$ws.Range("R4").Value = 'Región de O''Higgins'
$ws.Range("S4").Value = 1200

$ws.Range("D5").Value = 44698
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 220
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 15000
$ws.Range("Q5").Value = '$/caja 15 kilos granel'
$ws.Range("R5").Value = 'Región de O''Higgins'
$ws.Range("S5").Value = 1000

$ws.Range("D6").Value = 44698
$ws.Range("L6").Value = 'Segunda'
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 12000
$ws.Range("O6").Value = 12000
$ws.Range("P6").Value = 12000
$ws.Range("Q6").Value = '$/caja 15 kilos granel'
$ws.Range("R6").Value = 'Región de O''Higgins'
$ws.Range("S6").Value = 800

$ws.Range("D7").Value = 45068
$ws.Range("L7").Value = 'Primera'
$ws.Range("M7").Value = 350
$ws.Range("N7").Value = 10500
$ws.Range("O7").Value = 11000
$ws.Range("P7").Value = 10786
$ws.Range("Q7").Value = '$/caja 15 kilos granel'
$ws.Range("R7").Value = 'Provincia de Curicó'
$ws.Range("S7").Value = 719

$ws.Range("D8").Value = 44694
$ws.Range("M8").Value = 200
$ws.Range("N8").Value = 21600
$ws.Range("O8").Value = 21600
$ws.Range("P8").Value = 21600
$ws.Range("Q8").Value = '$/caja 18 kilos granel'
$ws.Range("R8").Value = 'Región de O''Higgins'
$ws.Range("S8").Value = 1200
$ws.Range("T8").Value = 18

$ws.Range("D9").Value = 44694
$ws.Range("N9").Value = 18000
$ws.Range("O9").Value = 18000
$ws.Range("P9").Value = 18000
$ws.Range("Q9").Value = '$/caja 18 kilos granel'
$ws.Range("R9").Value = 'Región de O''Higgins'
$ws.Range("S9").Value = 1000
$ws.Range("T9").Value = 18

$ws.Range("D10").Value = 44694
$ws.Range("L10").Value = 'Segunda'
$ws.Range("M10").Value = 250
$ws.Range("N10").Value = 14400
$ws.Range("O10").Value = 14400
$ws.Range("P10").Value = 14400
$ws.Range("Q10").Value = '$/caja 18 kilos granel'
$ws.Range("R10").Value = 'Región de O''Higgins'
$ws.Range("S10").Value = 800
$ws.Range("T10").Value = 18

$ws.Range("D11").Value = 45083
$ws.Range("L11").Value = 'Especial'
$ws.Range("M11").Value = 150
$ws.Range("R11").Value = 'Región Metropolitana'

$ws.Range("D12").Value = 45083
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 220
$ws.Range("N12").Value = 6000
$ws.Range("O12").Value = 6000
$ws.Range("P12").Value = 6000
$ws.Range("R12").Value = 'Región Metropolitana'
$ws.Range("S12").Value = 400

$ws.Range("D13").Value = 44309
$ws.Range("M13").Value = 40
$ws.Range("N13").Value = 18000
$ws.Range("O13").Value = 18000
$ws.Range("P13").Value = 18000
$ws.Range("S13").Value = 1200

$ws.Range("D14").Value = 44309
$ws.Range("M14").Value = 70
$ws.Range("N14").Value = 15000
$ws.Range("O14").Value = 15000
$ws.Range("P14").Value = 15000
$ws.Range("S14").Value = 1000

$ws.Range("D15").Value = 44685
$ws.Range("M15").Value = 350
$ws.Range("N15").Value = 21000
$ws.Range("O15").Value = 21000
$ws.Range("P15").Value = 21000
$ws.Range("Q15").Value = '$/caja 18 kilos granel'
$ws.Range("R15").Value = 'Región de O''Higgins'
$ws.Range("S15").Value = 1167
$ws.Range("T15").Value = 18

$ws.Range("D16").Value = 44685
$ws.Range("M16").Value = 330
$ws.Range("N16").Value = 15000
$ws.Range("O16").Value = 15000
$ws.Range("P16").Value = 15000
$ws.Range("Q16").Value = '$/caja 18 kilos granel'
$ws.Range("R16").Value = 'Región de O''Higgins'
$ws.Range("S16").Value = 833
$ws.Range("T16").Value = 18

$ws.Range("D17").Value = 44685
$ws.Range("L17").Value = 'Segunda'
$ws.Range("M17").Value = 280
$ws.Range("N17").Value = 10000
$ws.Range("O17").Value = 10000
$ws.Range("P17").Value = 10000
$ws.Range("Q17").Value = '$/caja 18 kilos granel'
$ws.Range("R17").Value = 'Región de O''Higgins'
$ws.Range("S17").Value = 556
$ws.Range("T17").Value = 18

$ws.Range("D18").Value = 45034
$ws.Range("M18").Value = 300
$ws.Range("N18").Value = 9000
$ws.Range("O18").Value = 9000
$ws.Range("P18").Value = 9000
$ws.Range("R18").Value = 'Paine'
$ws.Range("S18").Value = 600

$ws.Range("D19").Value = 45034
$ws.Range("L19").Value = 'Segunda'
$ws.Range("M19").Value = 280
$ws.Range("O19").Value = 6000
$ws.Range("P19").Value = 6000
$ws.Range("S19").Value = 400

$ws.Range("D20").Value = 45051
$ws.Range("L20").Value = 'Primera'
$ws.Range("N20").Value = 10500
$ws.Range("O20").Value = 10500
$ws.Range("P20").Value = 10500
$ws.Range("Q20").Value = '$/caja 15 kilos granel'
$ws.Range("R20").Value = 'Provincia de Curicó'
$ws.Range("S20").Value = 700
$ws.Range("T20").Value = 15

$ws.Range("D21").Value = 45051
$ws.Range("L21").Value = 'Segunda'
$ws.Range("M21").Value = 280
$ws.Range("N21").Value = 9000
$ws.Range("O21").Value = 9000
$ws.Range("P21").Value = 9000
$ws.Range("Q21").Value = '$/caja 15 kilos granel'
$ws.Range("R21").Value = 'Provincia de Curicó'
$ws.Range("S21").Value = 600
$ws.Range("T21").Value = 15

$ws.Range("D22").Value = 44285
$ws.Range("L22").Value = 'Especial'
$ws.Range("M22").Value = 40
$ws.Range("N22").Value = 18000
$ws.Range("O22").Value = 18000
$ws.Range("P22").Value = 18000
$ws.Range("Q22").Value = '$/caja 15 kilos empedrada'
$ws.Range("R22").Value = 'Provincia del Elquí'
$ws.Range("S22").Value = 1200
$ws.Range("T22").Value = 15

$ws.Range("D23").Value = 44285
$ws.Range("K23").Value = 'Wonderfull'
$ws.Range("L23").Value = 'Primera'
$ws.Range("M23").Value = 90
$ws.Range("Q23").Value = '$/caja 15 kilos empedrada'
$ws.Range("R23").Value = 'Provincia del Elquí'

$ws.Range("D24").Value = 44285
$ws.Range("K24").Value = 'Wonderfull'
$ws.Range("L24").Value = 'Segunda'
$ws.Range("M24").Value = 75
$ws.Range("Q24").Value = '$/caja 15 kilos empedrada'
$ws.Range("R24").Value = 'Provincia del Elquí'

$ws.Range("D25").Value = 44658
$ws.Range("K25").Value = 'Sin especificar'
$ws.Range("L25").Value = 'Especial'
$ws.Range("M25").Value = 280
$ws.Range("N25").Value = 21600
$ws.Range("O25").Value = 21600
$ws.Range("P25").Value = 21600
$ws.Range("Q25").Value = '$/caja 18 kilos granel'
$ws.Range("R25").Value = 'Provincia de Limarí'
$ws.Range("T25").Value = 18

$ws.Range("D26").Value = 44658
$ws.Range("K26").Value = 'Sin especificar'
$ws.Range("L26").Value = 'Primera'
$ws.Range("M26").Value = 330
$ws.Range("N26").Value = 16200
$ws.Range("O26").Value = 16200
$ws.Range("P26").Value = 16200
$ws.Range("Q26").Value = '$/caja 18 kilos granel'
$ws.Range("R26").Value = 'Provincia de Limarí'
$ws.Range("S26").Value = 900
$ws.Range("T26").Value = 18

$ws.Range("D27").Value = 44658
$ws.Range("K27").Value = 'Sin especificar'
$ws.Range("L27").Value = 'Segunda'
$ws.Range("M27").Value = 220
$ws.Range("N27").Value = 14400
$ws.Range("O27").Value = 14400
$ws.Range("P27").Value = 14400
$ws.Range("Q27").Value = '$/caja 18 kilos granel'
$ws.Range("R27").Value = 'Provincia de Limarí'
$ws.Range("S27").Value = 800
$ws.Range("T27").Value = 18

$ws.Range("D28").Value = 45063
$ws.Range("L28").Value = 'Especial'
$ws.Range("M28").Value = 220
$ws.Range("N28").Value = 10500
$ws.Range("O28").Value = 10500
$ws.Range("P28").Value = 10500
$ws.Range("R28").Value = 'Provincia de Curicó'
$ws.Range("S28").Value = 700

$ws.Range("D29").Value = 45063
$ws.Range("M29").Value = 250
$ws.Range("R29").Value = 'Provincia de Curicó'

$ws.Range("D30").Value = 44305
$ws.Range("L30").Value = 'Primera'
$ws.Range("M30").Value = 50
$ws.Range("N30").Value = 18000
$ws.Range("O30").Value = 18000
$ws.Range("P30").Value = 18000
$ws.Range("R30").Value = 'Región de O''Higgins'
$ws.Range("S30").Value = 1200

$ws.Range("D31").Value = 44305
$ws.Range("K31").Value = 'Wonderfull'
$ws.Range("L31").Value = 'Segunda'
$ws.Range("M31").Value = 60
$ws.Range("N31").Value = 15000
$ws.Range("O31").Value = 15000
$ws.Range("P31").Value = 15000
$ws.Range("Q31").Value = '$/caja 15 kilos granel'
$ws.Range("R31").Value = 'Región de O''Higgins'
$ws.Range("S31").Value = 1000
$ws.Range("T31").Value = 15

$ws.Range("D32").Value = 44687
$ws.Range("K32").Value = 'Wonderfull'
$ws.Range("L32").Value = 'Especial'
$ws.Range("M32").Value = 220
$ws.Range("N32").Value = 21000
$ws.Range("O32").Value = 21000
$ws.Range("P32").Value = 21000
$ws.Range("R32").Value = 'Región de O''Higgins'
$ws.Range("S32").Value = 1167

$ws.Range("D33").Value = 44687
$ws.Range("K33").Value = 'Wonderfull'
$ws.Range("L33").Value = 'Primera'
$ws.Range("M33").Value = 250
$ws.Range("N33").Value = 15000
$ws.Range("O33").Value = 15000
$ws.Range("P33").Value = 15000
$ws.Range("R33").Value = 'Región de O''Higgins'
$ws.Range("S33").Value = 833

$ws.Range("D34").Value = 44687
$ws.Range("L34").Value = 'Segunda'
$ws.Range("M34").Value = 280
$ws.Range("N34").Value = 10000
$ws.Range("O34").Value = 10000
$ws.Range("P34").Value = 10000
$ws.Range("S34").Value = 556

$ws.Range("D35").Value = 45076
$ws.Range("M35").Value = 470
$ws.Range("N35").Value = 7500
$ws.Range("O35").Value = 8000
$ws.Range("P35").Value = 7734
$ws.Range("Q35").Value = '$/caja 15 kilos granel'
$ws.Range("R35").Value = 'Región Metropolitana'
$ws.Range("S35").Value = 516
$ws.Range("T35").Value = 15

$ws.Range("D36").Value = 44678
$ws.Range("K36").Value = 'Sin especificar'
$ws.Range("L36").Value = 'Especial'
$ws.Range("M36").Value = 290
$ws.Range("N36").Value = 15000
$ws.Range("O36").Value = 15000
$ws.Range("P36").Value = 15000
$ws.Range("Q36").Value = '$/caja 15 kilos granel'
$ws.Range("S36").Value = 1000
$ws.Range("T36").Value = 15

$ws.Range("D37").Value = 44678
$ws.Range("K37").Value = 'Sin especificar'
$ws.Range("L37").Value = 'Primera'
$ws.Range("M37").Value = 220
$ws.Range("N37").Value = 12000
$ws.Range("O37").Value = 12000
$ws.Range("P37").Value = 12000
$ws.Range("Q37").Value = '$/caja 15 kilos granel'
$ws.Range("S37").Value = 800
$ws.Range("T37").Value = 15

$ws.Range("D38").Value = 44649
$ws.Range("K38").Value = 'Sin especificar'
$ws.Range("L38").Value = 'Especial'
$ws.Range("N38").Value = 21600
$ws.Range("O38").Value = 21600
$ws.Range("P38").Value = 21600
$ws.Range("R38").Value = 'Provincia de Limarí'
$ws.Range("S38").Value = 1200

$ws.Range("D39").Value = 44649
$ws.Range("K39").Value = 'Sin especificar'
$ws.Range("L39").Value = 'Primera'
$ws.Range("M39").Value = 250
$ws.Range("N39").Value = 16200
$ws.Range("O39").Value = 16200
$ws.Range("P39").Value = 16200
$ws.Range("R39").Value = 'Provincia de Limarí'
$ws.Range("S39").Value = 900

$ws.Range("D40").Value = 44649
$ws.Range("K40").Value = 'Sin especificar'
$ws.Range("L40").Value = 'Segunda'
$ws.Range("M40").Value = 180
$ws.Range("N40").Value = 14400
$ws.Range("O40").Value = 14400
$ws.Range("P40").Value = 14400
$ws.Range("R40").Value = 'Provincia de Limarí'
$ws.Range("S40").Value = 800

$ws.Range("D41").Value = 45085
$ws.Range("L41").Value = 'Primera'
$ws.Range("M41").Value = 400
$ws.Range("N41").Value = 6000
$ws.Range("O41").Value = 6500
$ws.Range("P41").Value = 6275
$ws.Range("Q41").Value = '$/caja 15 kilos granel'
$ws.Range("R41").Value = 'Paine'
$ws.Range("S41").Value = 418
$ws.Range("T41").Value = 15

$ws.Range("D42").Value = 45055
$ws.Range("K42").Value = 'Wonderfull'
$ws.Range("L42").Value = 'Primera'
$ws.Range("M42").Value = 470
$ws.Range("N42").Value = 10500
$ws.Range("O42").Value = 11000
$ws.Range("P42").Value = 10734
$ws.Range("Q42").Value = '$/caja 15 kilos granel'
$ws.Range("R42").Value = 'Provincia de Curicó'
$ws.Range("S42").Value = 716
$ws.Range("T42").Value = 15

$ws.Range("D43").Value = 45062
$ws.Range("K43").Value = 'Wonderfull'
$ws.Range("L43").Value = 'Especial'
$ws.Range("M43").Value = 200
$ws.Range("N43").Value = 10500
$ws.Range("O43").Value = 10500
$ws.Range("P43").Value = 10500
$ws.Range("Q43").Value = '$/caja 15 kilos granel'
$ws.Range("R43").Value = 'Provincia de Curicó'
$ws.Range("S43").Value = 700
$ws.Range("T43").Value = 15

$ws.Range("D44").Value = 45062
$ws.Range("K44").Value = 'Wonderfull'
$ws.Range("L44").Value = 'Primera'
$ws.Range("M44").Value = 200
$ws.Range("N44").Value = 9000
$ws.Range("O44").Value = 9000
$ws.Range("P44").Value = 9000
$ws.Range("Q44").Value = '$/caja 15 kilos granel'
$ws.Range("R44").Value = 'Provincia de Curicó'
$ws.Range("S44").Value = 600
$ws.Range("T44").Value = 15

$ws.Range("D45").Value = 45079
$ws.Range("K45").Value = 'Wonderfull'
$ws.Range("L45").Value = 'Primera'
$ws.Range("M45").Value = 200
$ws.Range("N45").Value = 7500
$ws.Range("O45").Value = 7500
$ws.Range("P45").Value = 7500
$ws.Range("Q45").Value = '$/caja 15 kilos granel'
$ws.Range("R45").Value = 'Provincia de Los Andes'
$ws.Range("S45").Value = 500
$ws.Range("T45").Value = 15

$ws.Range("D46").Value = 45070
$ws.Range("K46").Value = 'Wonderfull'
$ws.Range("L46").Value = 'Especial'
$ws.Range("M46").Value = 280
$ws.Range("N46").Value = 10500
$ws.Range("O46").Value = 10500
$ws.Range("P46").Value = 10500
$ws.Range("Q46").Value = '$/caja 15 kilos granel'
$ws.Range("R46").Value = 'Paine'
$ws.Range("S46").Value = 700
$ws.Range("T46").Value = 15

$ws.Range("D47").Value = 45070
$ws.Range("K47").Value = 'Wonderfull'
$ws.Range("L47").Value = 'Primera'
$ws.Range("N47").Value = 7500
$ws.Range("O47").Value = 7500
$ws.Range("P47").Value = 7500
$ws.Range("Q47").Value = '$/caja 15 kilos granel'
$ws.Range("R47").Value = 'Paine'
$ws.Range("S47").Value = 500
$ws.Range("T47").Value = 15

$ws.Range("D48").Value = 44664
$ws.Range("K48").Value = 'Sin especificar'
$ws.Range("L48").Value = 'Especial'
$ws.Range("M48").Value = 300
$ws.Range("N48").Value = 21600
$ws.Range("O48").Value = 21600
$ws.Range("P48").Value = 21600
$ws.Range("Q48").Value = '$/caja 18 kilos granel'
$ws.Range("R48").Value = 'Provincia de Limarí'
$ws.Range("T48").Value = 18

$ws.Range("D49").Value = 44664
$ws.Range("K49").Value = 'Sin especificar'
$ws.Range("L49").Value = 'Primera'
$ws.Range("M49").Value = 250
$ws.Range("N49").Value = 18000
$ws.Range("O49").Value = 18000
$ws.Range("P49").Value = 18000
$ws.Range("Q49").Value = '$/caja 18 kilos granel'
$ws.Range("R49").Value = 'Provincia de Limarí'
$ws.Range("T49").Value = 18

$ws.Range("D50").Value = 44664
$ws.Range("K50").Value = 'Sin especificar'
$ws.Range("L50").Value = 'Segunda'
$ws.Range("M50").Value = 250
$ws.Range("N50").Value = 16000
$ws.Range("O50").Value = 16000
$ws.Range("P50").Value = 16000
$ws.Range("Q50").Value = '$/caja 18 kilos granel'
$ws.Range("R50").Value = 'Provincia de Limarí'
$ws.Range("S50").Value = 889
$ws.Range("T50").Value = 18
